$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "A1 text: $($ws.Range("A1").Text)"
Write-Host "AZ1: $($ws.Range("AZ1").Text)"
Write-Host "BB2: $($ws.Range("BB2").Text)"
Write-Host "dim: $($ws.UsedRange.Address)"
